# Filter currency results endpoints added
# Updates the currency quote values and the report date/time text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated currency values (column B, rows 2-13)
$ws.Range("B2").Value = 5.17
$ws.Range("B3").Value = 3.62
$ws.Range("B4").Value = 4.04
$ws.Range("B5").Value = 5.29
$ws.Range("B6").Value = 5.44
$ws.Range("B7").Value = 0.0388
$ws.Range("B8").Value = 6.3
$ws.Range("B9").Value = 0.04
$ws.Range("B10").Value = 0.0057
$ws.Range("B12").Value = 0.25
$ws.Range("B13").Value = 0.77

# Updated report date and time (plain text cells)
$ws.Range("C16").Value = "31/07/2022"
$ws.Range("D16").Value = "22:46"
